# QE_holdings.xlsx update: refresh "as of" date and recompute Weight / Percent Change
# columns for each holding row (model holdings refresh from 2021-05-07 -> 2021-05-10).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet ships password-protected; unprotect so the cells can be written,
# then re-apply the same password once the data refresh is complete.
$ws.Unprotect("D382")

# --- Update the "as of" date in the confidentiality / disclosure footer ---
$footerCell = $ws.Range("A42")
$footerText = $footerCell.Value()
$footerText = $footerText -replace "2021-05-07", "2021-05-10"
$footerCell.Value = $footerText

# --- Refresh Weight (D) and Percent Change (E) for each holding row ---
$rows = @(
    @{ Row=2; D=0.05872000513469558; E=-0.02580446970278794 },
    @{ Row=3; D=0.05295364552897643; E=-0.02091420423037305 },
    @{ Row=4; D=0.309067561938314; E=-0.002375296912114022 },
    @{ Row=5; D=0.03452086452500081; E=-0.03072052885973742 },
    @{ Row=6; D=0.03164675643714214; E=-0.02649491642253998 },
    @{ Row=7; D=0.03043817327332228; E=-0.0001240387000744914 },
    @{ Row=8; D=0.02827438571380029; E=0.01050445103857567 },
    @{ Row=9; D=0.02352563131795134; E=0.004422253922967245 },
    @{ Row=10; D=0.02466594064979908; E=-0.02558749622650325 },
    @{ Row=11; D=0.02342453151146728; E=-0.04108687476494921 },
    @{ Row=12; D=0.02300293273339058; E=-0.003082029397818986 },
    @{ Row=13; D=0.02132366914270301; E=-0.002921445574550874 },
    @{ Row=14; D=0.02126084882726322; E=0.01864594894561589 },
    @{ Row=15; D=0.02134737096455508; E=0.005512159174649867 },
    @{ Row=16; D=0.02190058826663333; E=0.0077576918472404 },
    @{ Row=17; D=0.01968509717537209; E=-0.02413425679275449 },
    @{ Row=18; D=0.01451558300398988; E=-0.02947806485174265 },
    @{ Row=19; D=0.01706405322967737; E=-0.01032524522457401 },
    @{ Row=20; D=0.01539569667205446; E=0.01362397820163497 },
    @{ Row=21; D=0.01702315185568485; E=0.002402691013935421 },
    @{ Row=22; D=0.01410300350325512; E=-0.06444368428097635 },
    @{ Row=23; D=0.01486355930888005; E=0.007338103100348548 },
    @{ Row=24; D=0.01450299796583833; E=0.01461442786069678 },
    @{ Row=25; D=0.01400872059243648; E=-0.005053340819764163 },
    @{ Row=26; D=0.01372241097448887; E=-0.0000840689365280145 },
    @{ Row=27; D=0.01300758080748122; E=-0.01023147812205194 },
    @{ Row=28; D=0.01384605897432778; E=-0.004090165424468384 },
    @{ Row=29; D=0.01411328128441221; E=0.00707427993936327 },
    @{ Row=30; D=0.01349116089845427; E=-0.003109452736318463 },
    @{ Row=31; D=0.01221252102225748; E=0.007969222313822355 },
    @{ Row=32; D=0.01335324985537694; E=0.004060475161987043 },
    @{ Row=33; D=0.01233491051828123; E=-0.003060834077285968 },
    @{ Row=34; D=0.006213757712006504; E=-0.03689513747067463 },
    @{ Row=35; D=0.005284038018561254; E=-0.0340385836773579 },
    @{ Row=36; D=0.005314242110124957; E=-0.03840385222608156 },
    @{ Row=37; D=0.00512557141316974; E=-0.01913121764573489 },
    @{ Row=38; D=0.004746447138854502; E=-0.03641345176543365 },
    @{ Row=39; E=-0.008291757261461918 }
)

foreach ($row in $rows) {
    if ($row.ContainsKey("D")) {
        $ws.Cells.Item($row.Row, 4).Value = $row.D
    }
    if ($row.ContainsKey("E")) {
        $ws.Cells.Item($row.Row, 5).Value = $row.E
    }
}

# --- Restore sheet protection with the original password ---
$ws.Protect("D382")
